$wb = $excel.ActiveWorkbook

# --- Existing sheet: TestInputData ---
$ws1 = $wb.Worksheets.Item("TestInputData")

# Fix header/value placement + casing (A2 now "Chrome", B1 now "Product Name",
# B2 gets proper capitalisation "Samsung Galaxy S24 Ultra")
$ws1.Range("A2").Value = "Chrome"
$ws1.Range("B1").Value = "Product Name"
$ws1.Range("B2").Value = "Samsung Galaxy S24 Ultra"

# Widen column B slightly to fit the new text
$ws1.Columns.Item(2).ColumnWidth = 22

# Move the selection from D7 to B3
$null = $ws1.Range("B3").Select()

# --- New sheet: TestOutputData (added right after TestInputData) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TestOutputData"

$ws2.Range("A1").Value = "Product Name"
$ws2.Range("B1").Value = "Product Price"
$ws2.Range("A2").Value = "SAMSUNG Galaxy S24 Ultra 5G (Titanium Violet, 256 GB)"
$ws2.Range("B2").Value = "₹1,29,999"

$ws2.Columns.Item(1).ColumnWidth = 48
$ws2.Columns.Item(2).ColumnWidth = 11
